$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date cell number format (column D) before inserting, so we can re-apply it.
$dateNumberFormat = $ws.Cells.Item(18, 4).NumberFormat

# Insert a new row at position 18, shifting current row 18 (and below) down to row 19.
$ws.Rows.Item(18).Insert()

# New row 18: a new, more recent price record for "Especial"->"Primera" quality drop
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 45180
$ws.Cells.Item(18, 4).NumberFormat = $dateNumberFormat
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100107
$ws.Cells.Item(18, 8).Value = "Otros"
$ws.Cells.Item(18, 9).Value = 100107002
$ws.Cells.Item(18, 10).Value = "Chirimoya"
$ws.Cells.Item(18, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 40
$ws.Cells.Item(18, 14).Value = 22000
$ws.Cells.Item(18, 15).Value = 22000
$ws.Cells.Item(18, 16).Value = 22000
$ws.Cells.Item(18, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 2200
$ws.Cells.Item(18, 20).Value = 10

$wb.Save()
